$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.628571629524231
$ws.Range("B1").Value = 3.541851043701172
$ws.Range("C1").Value = 4.16127872467041
$ws.Range("D1").Value = 1.301671981811523
$ws.Range("E1").Value = 0.7635176181793213
